$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-11, columns A-G
# A: Colaborador_id, B: Colaborador_nome, C: Departamento,
# D: Motivo_da_ausência, E: Horas_de_ausência, F: Data_da_ausência, G: Salário

$data = @(
    @(83805, "Helena Campos",        "Atendimento ao Cliente", "Consulta médica",    4, 45105, 9870.26),
    @(54682, "Ana Sophia Almeida",   "P&D",                    "Doença",              2, 45082, 7760.95),
    @(50402, "Igor Martins",         "Vendas",                 "Problemas pessoais",  5, 45082, 7367.07),
    @(53351, "Vitor da Cunha",       "Recursos Humanos",        "Problemas pessoais",  6, 45097, 8413.370000000001),
    @(31054, "Manuela Cardoso",      "Engenharia",              "Outros",              4, 45095, 7515.04),
    @(55771, "Luiz Felipe da Mota",  "Atendimento ao Cliente", "Outros",              3, 45096, 7417.08),
    @(17830, "Sr. Nathan Cavalcanti","Engenharia",              "Problemas pessoais",  1, 45097, 9775.389999999999),
    @(79166, "Ana Sophia Ferreira",  "P&D",                    "Outros",              6, 45094, 9690.41),
    @(33916, "Benício Cunha",        "Recursos Humanos",        "Doença",              5, 45086, 6764.08),
    @(33406, "Vitor Hugo Farias",    "P&D",                    "Consulta médica",     4, 45086, 7849.65)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}
